# Append: 2025-09-25 01:43 JST
# - Update the "取得日時" (retrieved-at) timestamp in column A for every
#   data row (rows 2-23) from 2025-09-25 01:15:28 to 2025-09-25 01:43:36.
# - Two listings swapped positions on the source site between scrapes,
#   so the title (column B) and URL (column F) for rows 16 and 17 swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-25 01:43:36"

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Swap title/URL between row 16 and row 17
$b16 = $ws.Range("B16").Value()
$f16 = $ws.Range("F16").Value()
$b17 = $ws.Range("B17").Value()
$f17 = $ws.Range("F17").Value()

$ws.Range("B16").Value = $b17
$ws.Range("F16").Value = $f17
$ws.Range("B17").Value = $b16
$ws.Range("F17").Value = $f16
